$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.851.34"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.378.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.56%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.78"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +5.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.68"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -9.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.636"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.621"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.79"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -9.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0919"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.37"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.98%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.57%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.39"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.743.86"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.387.26"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.853.13"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +6.98%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.73"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +7.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.27"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "272.48"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +7.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.32"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -8.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.86"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +9.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.44"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.64"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.67%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.73"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0900"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.96"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -10.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.89"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.58"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -8.73%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.15%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -7.18%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.84"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.51"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.95"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +52.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.226"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "68.40"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.51%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "115.72"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.79"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.38%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.99"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.615.46"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +8.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.25"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.69%  "
